# Generate Report for Handoff
# Updates status text and handoff/handback timestamps across the
# Overview, zh-cn, and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status column (B2, C2) and Latest Handoff Date (D2)
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-23-11 12:23:20"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-11 12:23:17"

# de-de sheet: Status (C2) and Latest Handoff Datetime (E2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-11 12:23:20"
